# Lecture 4 and 5 Update
#
# 1) Refresh the cached "datetime1" field text on every Date Placeholder
#    shape across both Slide Masters and all Slide Layouts
#    (9/5/2018 -> 1/15/2019).
# 2) Update the title on the "Data Transformation II Info" slide to
#    just "Information".

$p = $ppt.ActivePresentation

function Update-DatePlaceholders($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "1/15/2019"
        }
    }
}

# --- Slide Masters (one per Design) ---
for ($d = 1; $d -le $p.Designs.Count; $d++) {
    $master = $p.Designs.Item($d).SlideMaster
    Update-DatePlaceholders $master.Shapes
}

# --- Slide Layouts (flattened across all masters/designs) ---
$flatMaster = $p.SlideMaster
for ($i = 1; $i -le $flatMaster.CustomLayouts.Count; $i++) {
    $layout = $flatMaster.CustomLayouts.Item($i)
    Update-DatePlaceholders $layout.Shapes
}

# --- Slide title text update ---
for ($k = 1; $k -le $p.Slides.Count; $k++) {
    $slide = $p.Slides.Item($k)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "Data Transformation II Info") {
                $shp.TextFrame.TextRange.Text = "Information"
            }
        }
    }
}
